$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.448.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "'2.633.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.64%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'595.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").Value = "'167.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -2.50%  "

$ws.Range("D9").Value = "'2.632.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "

$ws.Range("E10").Value = "  -3.02%  "

$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").Value = "'0.363"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").Value = "'5.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("D14").Value = "'27.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("D15").Value = "'3.112.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.62%  "

$ws.Range("E16").Value = "  -1.95%  "

$ws.Range("D17").Value = "'67.437.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.43%  "

$ws.Range("D18").Value = "'2.620.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.12%  "

$ws.Range("D19").Value = "'11.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.54%  "

$ws.Range("D20").Value = "'8.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.09%  "

$ws.Range("D21").Value = "'357.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.02%  "

$ws.Range("E22").Value = "  -1.95%  "

$ws.Range("E23").Value = "  -3.59%  "

$ws.Range("E24").Value = "  -4.94%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").Value = "'10.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.97%  "

$ws.Range("D27").Value = "'69.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.66%  "

$ws.Range("E28").Value = "  -1.91%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("E30").Value = "  -1.83%  "

$ws.Range("D31").Value = "'548.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.95%  "

$ws.Range("E32").Value = "  -1.43%  "

$ws.Range("E33").Value = "  -3.42%  "

$ws.Range("E34").Value = "  -2.34%  "

$ws.Range("D35").Value = "'0.136"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.28%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("E37").Value = "  -4.32%  "

$ws.Range("D38").Value = "'157.54"
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").Value = "'19.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.80%  "

$ws.Range("D40").Value = "'0.366"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.12%  "

$ws.Range("E41").Value = "  -1.10%  "

$ws.Range("D42").Value = "'18.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.91%  "

$ws.Range("E43").Value = "  -2.07%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("E45").Value = "  -4.37%  "

$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("D47").Value = "'152.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.66%  "

$ws.Range("E48").Value = "  -2.16%  "

$ws.Range("E49").Value = "  -1.67%  "

$ws.Range("D50").Value = "'1.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.62%  "

$ws.Range("E51").Value = "  -1.16%  "
